$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused helper columns F:G (previously H:I), shifting data left.
$ws.Range("F1:G1").EntireColumn.Delete()

# Update the intercept values (column E) for rows 2-7.
$ws.Range("E2").Value = -0.014704366287839665
$ws.Range("E3").Value = 0.15541563371216016
$ws.Range("E4").Value = 0.85837563371216019
$ws.Range("E5").Value = 0.69659563371216038
$ws.Range("E6").Value = 0.60597563371216001
$ws.Range("E7").Value = 0.77178563371216002

# Update the active selection.
$ws.Range("D6").Select() | Out-Null
